$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 100107.6
$ws.Range("J28").Value = 331627.66
$ws.Range("L28").Value = 331627.66
$ws.Range("N28").Value = -332597.66

$ws.Range("H88").Value = 1163.8
$ws.Range("I88").Value = 734
$ws.Range("J88").Value = 1450.3334
$ws.Range("K88").Value = 734
$ws.Range("L88").Value = 1450.3334
$ws.Range("M88").Value = -328
$ws.Range("N88").Value = -2262.3334

$ws.Range("H91").Value = 1163.8
$ws.Range("I91").Value = 734
$ws.Range("J91").Value = 1450.3334
$ws.Range("K91").Value = 734
$ws.Range("L91").Value = 1450.3334
$ws.Range("M91").Value = 670
$ws.Range("N91").Value = -4258.3334

$ws.Range("H113").Value = 7999
$ws.Range("I113").Value = 7998.25
$ws.Range("K113").Value = 7998.25
$ws.Range("M113").Value = -4744.25

$ws.Range("H132").Value = 1664.7059
$ws.Range("I132").Value = 1456.3125
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 4368.9375
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -1838.9375
$ws.Range("N132").Value = -20057

$ws.Range("H137").Value = 7004.2036
$ws.Range("I137").Value = 2734.2122
$ws.Range("J137").Value = 12423.808
$ws.Range("K137").Value = 8202.6366
$ws.Range("L137").Value = 37271.424
$ws.Range("M137").Value = -5652.6366
$ws.Range("N137").Value = -42371.424

$ws.Range("H138").Value = 4405.4326
$ws.Range("J138").Value = 5154.1304
$ws.Range("L138").Value = 15462.3912
$ws.Range("N138").Value = -25742.3912

$ws.Range("H140").Value = 156250
$ws.Range("I140").Value = 50000
$ws.Range("K140").Value = 50000
$ws.Range("M140").Value = -44820

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H40").Value = 15694.25
$ws.Range("I40").Value = 14851.777
$ws.Range("J40").Value = 18221.666
$ws.Range("K40").Value = 14851.777
$ws.Range("L40").Value = 18221.666
$ws.Range("M40").Value = -14675.777
$ws.Range("N40").Value = -18573.666

$ws.Range("H45").Value = 9668.0625
$ws.Range("I45").Value = 11437.615
$ws.Range("K45").Value = 11437.615
$ws.Range("M45").Value = -11060.615

$ws.Range("H61").Value = 12374.286
$ws.Range("I61").Value = 5000
$ws.Range("J61").Value = 17905
$ws.Range("K61").Value = 5000
$ws.Range("L61").Value = 17905
$ws.Range("M61").Value = -4788
$ws.Range("N61").Value = -18329

$ws.Range("H74").Value = 11674.117
$ws.Range("I74").Value = 11651
$ws.Range("J74").Value = 11749.25
$ws.Range("K74").Value = 11651
$ws.Range("L74").Value = 11749.25
$ws.Range("M74").Value = -10777
$ws.Range("N74").Value = -13497.25

$ws.Range("H77").Value = 11674.117
$ws.Range("I77").Value = 11651
$ws.Range("J77").Value = 11749.25
$ws.Range("K77").Value = 58255
$ws.Range("L77").Value = 58746.25
$ws.Range("M77").Value = -53887
$ws.Range("N77").Value = -67482.25

$ws.Range("H92").Value = 20550
$ws.Range("J92").Value = 20550
$ws.Range("L92").Value = 20550
$ws.Range("N92").Value = -25542

$ws.Range("H132").Value = 4032.3062
$ws.Range("I132").Value = 2940.1
$ws.Range("J132").Value = 8886.556
$ws.Range("K132").Value = 8820.299999999999
$ws.Range("L132").Value = 26659.668
$ws.Range("M132").Value = -6290.299999999999
$ws.Range("N132").Value = -31719.668

$ws.Range("H136").Value = 12374.286
$ws.Range("I136").Value = 5000
$ws.Range("J136").Value = 17905
$ws.Range("K136").Value = 15000
$ws.Range("L136").Value = 53715
$ws.Range("M136").Value = -12450
$ws.Range("N136").Value = -58815

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3704.5
$ws.Range("I20").Value = 2624.2727
$ws.Range("K20").Value = 2624.2727
$ws.Range("M20").Value = -2377.2727

$ws.Range("H107").Value = 814.0833
$ws.Range("I107").Value = 835.44446
$ws.Range("K107").Value = 835.44446
$ws.Range("M107").Value = 1084.55554

$ws.Range("H134").Value = 7376.5903
$ws.Range("I134").Value = 6569.909
$ws.Range("J134").Value = 9464.471
$ws.Range("K134").Value = 19709.727
$ws.Range("L134").Value = 28393.413
$ws.Range("M134").Value = -17174.727
$ws.Range("N134").Value = -33463.413

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4081.25
$ws.Range("I16").Value = 3719
$ws.Range("K16").Value = 3719
$ws.Range("M16").Value = -3432

$ws.Range("H58").Value = 28021.824
$ws.Range("I58").Value = 41778
$ws.Range("K58").Value = 41778
$ws.Range("M58").Value = -41575

$ws.Range("H62").Value = 238123
$ws.Range("I62").Value = 337673
$ws.Range("J62").Value = 208258
$ws.Range("K62").Value = 337673
$ws.Range("L62").Value = 208258
$ws.Range("M62").Value = -337049
$ws.Range("N62").Value = -209506

$ws.Range("H65").Value = 238123
$ws.Range("I65").Value = 337673
$ws.Range("J65").Value = 208258
$ws.Range("K65").Value = 1688365
$ws.Range("L65").Value = 1041290
$ws.Range("M65").Value = -1685245
$ws.Range("N65").Value = -1047530

$ws.Range("H105").Value = 2236.682
$ws.Range("I105").Value = 2192
$ws.Range("K105").Value = 2192
$ws.Range("M105").Value = -445

$ws.Range("H113").Value = 4081.25
$ws.Range("I113").Value = 3719
$ws.Range("K113").Value = 3719
$ws.Range("M113").Value = -1549

$ws.Range("H134").Value = 3550.2
$ws.Range("I134").Value = 1544.7059
$ws.Range("J134").Value = 7811.875
$ws.Range("K134").Value = 4634.1177
$ws.Range("L134").Value = 23435.625
$ws.Range("M134").Value = -2099.1177
$ws.Range("N134").Value = -28505.625

$ws.Range("H136").Value = 28021.824
$ws.Range("I136").Value = 41778
$ws.Range("K136").Value = 125334
$ws.Range("M136").Value = -122784

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 7640
$ws.Range("J39").Value = 8798
$ws.Range("L39").Value = 26394
$ws.Range("N39").Value = -26982

$ws.Range("H55").Value = 4977.7144
$ws.Range("J55").Value = 7238.25
$ws.Range("L55").Value = 21714.75
$ws.Range("N55").Value = -22068.75

$ws.Range("H119").Value = 1299.6666
$ws.Range("I119").Value = 1299.6666
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 3898.9998
$ws.Range("L119").Value = 0
$ws.Range("M119").Value = 939.0001999999999
$ws.Range("N119").ClearContents()

$ws.Range("H120").Value = 3000
$ws.Range("I120").Value = 3000
$ws.Range("K120").Value = 9000
$ws.Range("M120").Value = -4162

$ws.Range("H131").Value = 2100.075
$ws.Range("J131").Value = 2443.75
$ws.Range("L131").Value = 7331.25
$ws.Range("N131").Value = -17411.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2312.5
$ws.Range("J102").Value = 2779.5
$ws.Range("L102").Value = 2779.5
$ws.Range("N102").Value = -6023.5

$ws.Range("H113").Value = 190082.75
$ws.Range("I113").Value = 283187
$ws.Range("K113").Value = 283187
$ws.Range("M113").Value = -281017

$ws.Range("H126").Value = 3321.3333
$ws.Range("I126").Value = 3321.3333
$ws.Range("K126").Value = 9963.999899999999
$ws.Range("M126").Value = -7493.999899999999

$ws.Range("H132").Value = 15262.632
$ws.Range("I132").Value = 18155.428
$ws.Range("J132").Value = 7162.8
$ws.Range("K132").Value = 54466.284
$ws.Range("L132").Value = 21488.4
$ws.Range("M132").Value = -51936.284
$ws.Range("N132").Value = -26548.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5749.375
$ws.Range("I7").Value = 6799
$ws.Range("J7").Value = 4000
$ws.Range("K7").Value = 6799
$ws.Range("L7").Value = 4000
$ws.Range("M7").Value = -6687
$ws.Range("N7").Value = -4224

$ws.Range("H93").Value = 4829.3887
$ws.Range("I93").Value = 4671.154
$ws.Range("J93").Value = 5240.8
$ws.Range("K93").Value = 4671.154
$ws.Range("L93").Value = 5240.8
$ws.Range("M93").Value = -3423.154
$ws.Range("N93").Value = -7736.8

$ws.Range("H126").Value = 5749.375
$ws.Range("I126").Value = 6799
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 20397
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -17927
$ws.Range("N126").Value = -16940

$ws.Range("H136").Value = 8403
$ws.Range("I136").Value = 8256
$ws.Range("K136").Value = 24768
$ws.Range("M136").Value = -22218

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 24026
$ws.Range("I34").Value = 24026
$ws.Range("K34").Value = 24026
$ws.Range("M34").Value = -23823

$ws.Range("H37").Value = 46599.8
$ws.Range("J37").Value = 43000
$ws.Range("L37").Value = 43000
$ws.Range("N37").Value = -43406

$ws.Range("H40").Value = 19000
$ws.Range("J40").Value = 19000
$ws.Range("L40").Value = 19000
$ws.Range("N40").Value = -19298

$ws.Range("H42").Value = 70287.8
$ws.Range("I42").Value = 48199.5
$ws.Range("J42").Value = 85013.336
$ws.Range("K42").Value = 48199.5
$ws.Range("L42").Value = 85013.336
$ws.Range("M42").Value = -47821.5
$ws.Range("N42").Value = -85769.336

$ws.Range("H43").Value = 62250
$ws.Range("I43").Value = 49500
$ws.Range("K43").Value = 49500
$ws.Range("M43").Value = -49351

$ws.Range("H132").Value = 7623.2705
$ws.Range("I132").Value = 4005.7188
$ws.Range("J132").Value = 30775.6
$ws.Range("K132").Value = 12017.1564
$ws.Range("L132").Value = 92326.79999999999
$ws.Range("M132").Value = -9487.1564
$ws.Range("N132").Value = -97386.79999999999

$ws.Range("H136").Value = 2057.5454
$ws.Range("I136").Value = 1045.6428
$ws.Range("K136").Value = 3136.9284
$ws.Range("M136").Value = -586.9284000000002

$ws.Range("H140").Value = 34485.8
$ws.Range("J140").Value = 34485.8
$ws.Range("L140").Value = 34485.8
$ws.Range("N140").Value = -44845.8
